$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.151.93"
$ws.Range("E2").Value = "  -3.12%  "
$ws.Range("D3").Value = "1.652.25"
$ws.Range("E3").Value = "  -5.11%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'239.62"
$ws.Range("E5").Value = "  -3.58%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "'0.4796"
$ws.Range("E7").Value = "  -6.75%  "
$ws.Range("D8").Value = "'0.2641"
$ws.Range("E8").Value = "  -4.14%  "
$ws.Range("D9").Value = "'0.06012"
$ws.Range("E9").Value = "  -2.85%  "
$ws.Range("D10").Value = "'0.07152"
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("D11").Value = "1.649.89"
$ws.Range("E11").Value = "  -5.33%  "
$ws.Range("D12").Value = "'14.55"
$ws.Range("E12").Value = "  -3.39%  "
$ws.Range("D13").Value = "'0.6244"
$ws.Range("E13").Value = "  -3.65%  "
$ws.Range("D14").Value = "'4.604"
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("D15").Value = "'73.57"
$ws.Range("E15").Value = "  -5.17%  "
$ws.Range("D16").Value = "'0.9996"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("D18").Value = "25.152.44"
$ws.Range("E18").Value = "  -3.28%  "
$ws.Range("D19").Value = "'11.42"
$ws.Range("E19").Value = "  -2.97%  "
$ws.Range("D20").Value = "'0.000006585"
$ws.Range("E20").Value = "  -3.14%  "
$ws.Range("D21").Value = "'4.483"
$ws.Range("E21").Value = "  +5.34%  "
$ws.Range("D22").Value = "1.862.00"
$ws.Range("E22").Value = "  -5.49%  "
$ws.Range("D23").Value = "'8.648"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'5.321"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "'132.90"
$ws.Range("E25").Value = "  -2.47%  "
$ws.Range("D26").Value = "'14.84"
$ws.Range("E26").Value = "  -3.01%  "
$ws.Range("D27").Value = "'1.402"
$ws.Range("E27").Value = "  -7.25%  "
$ws.Range("B28").Value = "BitcoinCash"
$ws.Range("C28").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D28").Value = "'102.97"
$ws.Range("E28").Value = "  -2.48%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "'1.674"
$ws.Range("E29").Value = "  -6.00%  "
$ws.Range("E30").Value = "  -3.76%  "
$ws.Range("D31").Value = "'0.07935"
$ws.Range("E31").Value = "  -4.00%  "
$ws.Range("D32").Value = "'3.606"
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("D33").Value = "'0.04610"
$ws.Range("E33").Value = "  -1.63%  "
$ws.Range("D34").Value = "'2.626"
$ws.Range("E34").Value = "  -0.64%  "
$ws.Range("D35").Value = "'0.9457"
$ws.Range("E35").Value = "  -5.50%  "
$ws.Range("D36").Value = "'0.5818"
$ws.Range("E36").Value = "  -6.45%  "
$ws.Range("D37").Value = "'2.638"
$ws.Range("E37").Value = "  -3.06%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'0.8661"
$ws.Range("E38").Value = "  +14.85%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01556"
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("D40").Value = "'1.000"
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'1.834"
$ws.Range("E41").Value = "  -5.26%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "'99.31"
$ws.Range("E42").Value = "  -1.01%  "
$ws.Range("D43").Value = "'0.3710"
$ws.Range("E43").Value = "  -3.75%  "
$ws.Range("D44").Value = "'4.807"
$ws.Range("E44").Value = "  -4.03%  "
$ws.Range("D45").Value = "'0.1141"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "'6.070"
$ws.Range("E46").Value = "  -4.34%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.05185"
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("D48").Value = "'52.37"
$ws.Range("E48").Value = "  -5.16%  "
$ws.Range("D49").Value = "'29.78"
$ws.Range("E49").Value = "  -2.57%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.340"
$ws.Range("E51").Value = "  -4.34%  "
